$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Waybill Number"
$ws.Range("B1").Value = "WONumber"
$ws.Range("C1").Value = "ReferenceNumber"

# Data rows - Turkish + Nippon RPA refreshed waybill data
$ws.Range("A2").Value = "235-36155346"
$ws.Range("B2").Value = "DJSJUA4241863"
$ws.Range("C2").Value = 2340126729

$ws.Range("A3").Value = "235-36298102"
$ws.Range("B3").Value = "DJCVGA4241735"
$ws.Range("C3").Value = 2200250137

$ws.Range("A4").Value = "235-36298065"
$ws.Range("B4").Value = "DJCVGA4241725"
$ws.Range("C4").Value = 2200250134

$ws.Range("A5").Value = "235-36297586"
$ws.Range("B5").Value = "DJJNBA4241335"
$ws.Range("C5").Value = 2780468943

$ws.Range("A6").Value = "235-36155162"
$ws.Range("B6").Value = "DJSAVA4241086"
$ws.Range("C6").Value = 2330055233

$ws.Range("A7").Value = "235-36297586"
$ws.Range("B7").Value = "DJDURA4240721"
$ws.Range("C7").Value = 2790022451

$ws.Range("A8").Value = "235-33973601"
$ws.Range("B8").Value = "DJAMSA4239118"
$ws.Range("C8").Value = 2482459020

$ws.Range("A9").Value = "235-36297726"
$ws.Range("B9").Value = "DJSTRA4238475"
$ws.Range("C9").Value = 2570304342

$ws.Range("A10").Value = "235-39888052"
$ws.Range("B10").Value = "DJSINA4238206"
$ws.Range("C10").Value = 2711779354

$ws.Range("A11").Value = "235-36925383"
$ws.Range("B11").Value = "DJISTA4238100"
$ws.Range("C11").Value = 2640264509

# Update the sheet selection to match the refreshed data range
$ws.Range("A2:C11").Select()
